# Fruta / hortaliza, semanal
# Insert a new row of data at row 732, pushing existing rows 732:825
# down to 733:826 (dimension grows from A1:T825 to A1:T826).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 732; this shifts rows 732-825 down to 733-826.
$ws.Rows.Item(732).Insert()

# Populate the newly inserted row 732 with the new record's data.
$ws.Cells.Item(732, 1).Value = 10
$ws.Cells.Item(732, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(732, 3).Value = "La Araucanía"
$ws.Cells.Item(732, 4).Value = 45127
$ws.Cells.Item(732, 5).Value = 9
$ws.Cells.Item(732, 6).Value = "Fruta"
$ws.Cells.Item(732, 7).Value = 100101
$ws.Cells.Item(732, 8).Value = "Berries"
$ws.Cells.Item(732, 9).Value = 100101007
$ws.Cells.Item(732, 10).Value = "Kiwi"
$ws.Cells.Item(732, 11).Value = "Sin especificar"
$ws.Cells.Item(732, 12).Value = "Especial"
$ws.Cells.Item(732, 13).Value = 310
$ws.Cells.Item(732, 14).Value = 15000
$ws.Cells.Item(732, 15).Value = 15000
$ws.Cells.Item(732, 16).Value = 15000
$ws.Cells.Item(732, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(732, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(732, 19).Value = 1500
$ws.Cells.Item(732, 20).Value = 10
